$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Лист1")

# New rows describing additional server operations (Client-Server, server side)
# Row 28: Add floor to active ones
$ws.Range("A28").Value = "Добавить этаж к действующим"
$ws.Range("B28").Value = 24
$ws.Range("B28").HorizontalAlignment = -4108
$ws.Range("C28").Value = "id"
$ws.Range("C28").HorizontalAlignment = -4108
$ws.Range("D28").Value = "true/false"

# Row 29: Remove from active ones
$ws.Range("A29").Value = "Убрать из действующих"
$ws.Range("B29").Value = 25
$ws.Range("B29").HorizontalAlignment = -4108
$ws.Range("C29").Value = "id"
$ws.Range("C29").HorizontalAlignment = -4108
$ws.Range("D29").Value = "true/false"

# Row 30: Delete manager
$ws.Range("A30").Value = "Удалить менеджера"
$ws.Range("B30").Value = 26
$ws.Range("B30").HorizontalAlignment = -4108
$ws.Range("C30").Value = "id"
$ws.Range("C30").HorizontalAlignment = -4108

# Row 31: Change password
$ws.Range("A31").Value = "Сменить пароль"
$ws.Range("B31").Value = 27
$ws.Range("B31").HorizontalAlignment = -4108
$ws.Range("C31").Value = "login pass"
$ws.Range("C31").HorizontalAlignment = -4108

# Update view state: scroll down a bit and move the selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C32").Select()
